$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Corrected data-cleaning pass: the pre/post/total fixation metrics
#    (rows 3-8) had been computed against the wrong denominator; update
#    every affected cell to the corrected value.
# ------------------------------------------------------------------

# Row 3 - Revisit count
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 22
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 6
$ws.Range("W3").Value = 16
$ws.Range("Y3").Value = 41

# Row 4 - Fixation count
$ws.Range("G4").Value = 202
$ws.Range("H4").Value = 106
$ws.Range("I4").Value = 25
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 12
$ws.Range("L4").Value = 52
$ws.Range("M4").Value = 24
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 10
$ws.Range("W4").Value = 31
$ws.Range("Y4").Value = 161

# Row 5 - Dwell time (ms)
$ws.Range("G5").Value = 76278.87
$ws.Range("H5").Value = 44428.9
$ws.Range("I5").Value = 13480.06
$ws.Range("J5").Value = 8642.61
$ws.Range("K5").Value = 5223.06
$ws.Range("L5").Value = 21137.4
$ws.Range("M5").Value = 12745.94
$ws.Range("N5").Value = 6090.86
$ws.Range("O5").Value = 4221.19
$ws.Range("W5").Value = 11562.39
$ws.Range("Y5").Value = 67154.2

# Row 6 - Dwell time (%)
$ws.Range("D6").Value = 0.07
$ws.Range("E6").Value = 3.81
$ws.Range("F6").Value = 3.28
$ws.Range("G6").Value = 33.59
$ws.Range("H6").Value = 19.57
$ws.Range("I6").Value = 5.94
$ws.Range("J6").Value = 3.81
$ws.Range("K6").Value = 2.3
$ws.Range("L6").Value = 9.31
$ws.Range("M6").Value = 5.61
$ws.Range("N6").Value = 2.68
$ws.Range("O6").Value = 1.86
$ws.Range("Q6").Value = 1.04
$ws.Range("R6").Value = 0.21
$ws.Range("S6").Value = 0.77
$ws.Range("U6").Value = 0.28
$ws.Range("V6").Value = 0.43
$ws.Range("W6").Value = 5.09
$ws.Range("X6").Value = 0.44
$ws.Range("Y6").Value = 29.58
$ws.Range("Z6").Value = 2.65
$ws.Range("AA6").Value = 2.13

# Row 7 - Fixation duration (ms)
$ws.Range("G7").Value = 377.62
$ws.Range("H7").Value = 419.14
$ws.Range("I7").Value = 539.2
$ws.Range("J7").Value = 540.16
$ws.Range("K7").Value = 435.26
$ws.Range("L7").Value = 406.49
$ws.Range("M7").Value = 531.08
$ws.Range("N7").Value = 609.09
$ws.Range("O7").Value = 422.12
$ws.Range("W7").Value = 372.98
$ws.Range("Y7").Value = 417.11

# Row 8 - First fixation duration (ms)
$ws.Range("I8").Value = 116.69
$ws.Range("L8").Value = 150.14
$ws.Range("M8").Value = 116.69

# ------------------------------------------------------------------
# 2. The header row (A1:AA1) no longer carries the bold / bordered /
#    centred "table header" look - and the A1 label itself is blanked
#    out. Strip the formatting back to the workbook default and clear
#    the stray "Unnamed: 0" text.
# ------------------------------------------------------------------
$header = $ws.Range("A1:AA1")
$header.ClearFormats()
$ws.Range("A1").Value = ""

# ------------------------------------------------------------------
# 3. Two trailing, entirely-blank rows (10 and 11) were removed from
#    the sheet, shrinking the used range from A1:AA11 down to A1:AA9.
# ------------------------------------------------------------------
$ws.Range("A10:AA11").EntireRow.Delete()
